# This script reorders the player roster rows (A2:C18) on the active sheet.
# The underlying (player, position, team) triples are unchanged - only the
# row order changes: the block of 7 rows (Brandon Miller, Bilal Coulibaly,
# Cameron Johnson, Julius Randle, Anthony Davis, Cade Cunningham, Derrick
# White) is moved up to the top of the list (right after the header row),
# the remaining rows keep their relative order, and the last two rows
# (Brandon Ingram, LaMelo Ball) stay fixed at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Brandon Miller", "SG,SF", "Charlotte Hornets"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("Julius Randle", "PF", "Minnesota Timberwolves"),
    @("Anthony Davis", "PF,C", "Los Angeles Lakers"),
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Malik Monk", "SG,SF", "Sacramento Kings"),
    @("Herbert Jones", "SF,PF", "New Orleans Pelicans"),
    @("Duncan Robinson", "SG,SF", "Miami Heat"),
    @("Kelly Olynyk", "C", "Toronto Raptors"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Kentavious Caldwell-Pope", "SG,SF", "Orlando Magic"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}
